$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price (column D) and 1h volume-change
# (column E) figures to the latest scrape.
#
# Every cell in D/E is stored as literal text in the source workbook
# (coinranking prices like "56.473.18" use '.' as a thousands
# separator, and the percentages keep their padding spaces). Most new
# values are non-numeric-looking already (multi-dot prices, the
# padded "  +x.xx%  " strings, …) so a plain .Value assignment keeps
# them as text. A handful of new prices look like plain decimals
# (e.g. "19.00", "0.169") and Excel would silently reinterpret those
# as numbers -- losing the trailing zero / turning them into a Number
# cell -- so those are assigned with a leading apostrophe to force
# text, exactly as typing them into Excel by hand would.

$ws.Range("D2").Value = '56.473.18'
$ws.Range("E2").Value = '  +3.89%  '
$ws.Range("D3").Value = '2.988.75'
$ws.Range("E3").Value = '  +4.13%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''506.25'
$ws.Range("E5").Value = '  +8.32%  '
$ws.Range("D6").Value = '''137.83'
$ws.Range("E6").Value = '  +10.37%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +6.90%  '
$ws.Range("D9").Value = '''7.54'
$ws.Range("E9").Value = '  +14.64%  '
$ws.Range("E10").Value = '  +12.14%  '
$ws.Range("E11").Value = '  +5.68%  '
$ws.Range("E12").Value = '  +5.56%  '
$ws.Range("D13").Value = '3.502.23'
$ws.Range("E13").Value = '  +4.33%  '
$ws.Range("D14").Value = '''25.37'
$ws.Range("E14").Value = '  +9.37%  '
$ws.Range("E15").Value = '  +15.28%  '
$ws.Range("D16").Value = '56.524.89'
$ws.Range("E16").Value = '  +4.14%  '
$ws.Range("D17").Value = '2.990.37'
$ws.Range("E17").Value = '  +4.84%  '
$ws.Range("D18").Value = '''5.88'
$ws.Range("E18").Value = '  +10.11%  '
$ws.Range("D19").Value = '''12.35'
$ws.Range("E19").Value = '  +8.70%  '
$ws.Range("D20").Value = '''7.79'
$ws.Range("E20").Value = '  +11.14%  '
$ws.Range("D21").Value = '''325.92'
$ws.Range("E21").Value = '  +9.59%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("E23").Value = '  +9.28%  '
$ws.Range("D24").Value = '''62.31'
$ws.Range("E24").Value = '  +7.02%  '
$ws.Range("D25").Value = '''0.169'
$ws.Range("E25").Value = '  +12.69%  '
$ws.Range("E26").Value = '  +0.92%  '
$ws.Range("D27").Value = '0.0₃0902'
$ws.Range("E27").Value = '  +14.32%  '
$ws.Range("D28").Value = '''6.56'
$ws.Range("E28").Value = '  +7.77%  '
$ws.Range("D29").Value = '''7.04'
$ws.Range("E30").Value = '  +14.40%  '
$ws.Range("E31").Value = '  +10.99%  '
$ws.Range("D32").Value = '''20.53'
$ws.Range("E32").Value = '  +9.91%  '
$ws.Range("D33").Value = '''155.79'
$ws.Range("E33").Value = '  +10.55%  '
$ws.Range("E34").Value = '  +8.26%  '
$ws.Range("D35").Value = '''5.59'
$ws.Range("E35").Value = '  +4.30%  '
$ws.Range("E36").Value = '  +4.31%  '
$ws.Range("D37").Value = '''0.0676'
$ws.Range("E37").Value = '  +10.21%  '
$ws.Range("D38").Value = '''23.93'
$ws.Range("E38").Value = '  +4.65%  '
$ws.Range("D39").Value = '3.023.19'
$ws.Range("E39").Value = '  +4.43%  '
$ws.Range("E40").Value = '  +5.27%  '
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("E42").Value = '  +7.83%  '
$ws.Range("D43").Value = '2.253.34'
$ws.Range("E43").Value = '  +11.08%  '
$ws.Range("E44").Value = '  +7.96%  '
$ws.Range("D45").Value = '''0.982'
$ws.Range("E45").Value = '  +5.82%  '
$ws.Range("E46").Value = '  +5.86%  '
$ws.Range("D47").Value = '''1.99'
$ws.Range("E47").Value = '  +26.85%  '
$ws.Range("E48").Value = '  +10.67%  '
$ws.Range("D49").Value = '''5.75'
$ws.Range("E49").Value = '  +8.06%  '
$ws.Range("D50").Value = '''19.00'
$ws.Range("E50").Value = '  +8.24%  '
$ws.Range("E51").Value = '  +10.48%  '
